{"js": "// Replace the computed three-digit \u00f7 one-digit division answers in the\n// worksheet table with a newly generated set of problems/answers.\n// Each old answer string is unique in the document, so a straightforward\n// search-and-replace (by exact text) is safe and order-independent.\nconst replacements = [\n  [\"904\u00f73=301, 1\", \"241\u00f76=40, 1\"],\n  [\"882\u00f72=441, 0\", \"185\u00f76=30, 5\"],\n  [\"425\u00f76=70, 5\", \"268\u00f79=29, 7\"],\n  [\"364\u00f73=121, 1\", \"953\u00f77=136, 1\"],\n  [\"411\u00f79=45, 6\", \"715\u00f75=143, 0\"],\n  [\"912\u00f79=101, 3\", \"749\u00f73=249, 2\"],\n  [\"289\u00f77=41, 2\", \"787\u00f76=131, 1\"],\n  [\"906\u00f77=129, 3\", \"481\u00f72=240, 1\"],\n  [\"698\u00f79=77, 5\", \"689\u00f75=137, 4\"],\n  [\"862\u00f78=107, 6\", \"531\u00f73=177, 0\"],\n  [\"104\u00f76=17, 2\", \"803\u00f78=100, 3\"],\n  [\"653\u00f73=217, 2\", \"607\u00f73=202, 1\"],\n  [\"335\u00f76=55, 5\", \"635\u00f78=79, 3\"],\n  [\"615\u00f77=87, 6\", \"407\u00f73=135, 2\"],\n  [\"618\u00f77=88, 2\", \"305\u00f79=33, 8\"],\n  [\"241\u00f74=60, 1\", \"650\u00f77=92, 6\"],\n  [\"502\u00f78=62, 6\", \"370\u00f73=123, 1\"],\n  [\"279\u00f77=39, 6\", \"319\u00f78=39, 7\"],\n  [\"101\u00f75=20, 1\", \"311\u00f76=51, 5\"],\n  [\"701\u00f78=87, 5\", \"978\u00f78=122, 2\"],\n  [\"703\u00f72=351, 1\", \"286\u00f76=47, 4\"],\n  [\"974\u00f77=139, 1\", \"642\u00f77=91, 5\"],\n  [\"956\u00f78=119, 4\", \"333\u00f74=83, 1\"],\n  [\"984\u00f75=196, 4\", \"145\u00f78=18, 1\"],\n  [\"206\u00f73=68, 2\", \"896\u00f75=179, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the computed three-digit \u00f7 one-digit division answers in the\n# worksheet table with a newly generated set of problems/answers.\n# Each old answer string is unique in the document, so Find/Replace across\n# the whole document body is safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"904\u00f73=301, 1\"; New = \"241\u00f76=40, 1\" },\n    @{ Old = \"882\u00f72=441, 0\"; New = \"185\u00f76=30, 5\" },\n    @{ Old = \"425\u00f76=70, 5\"; New = \"268\u00f79=29, 7\" },\n    @{ Old = \"364\u00f73=121, 1\"; New = \"953\u00f77=136, 1\" },\n    @{ Old = \"411\u00f79=45, 6\"; New = \"715\u00f75=143, 0\" },\n    @{ Old = \"912\u00f79=101, 3\"; New = \"749\u00f73=249, 2\" },\n    @{ Old = \"289\u00f77=41, 2\"; New = \"787\u00f76=131, 1\" },\n    @{ Old = \"906\u00f77=129, 3\"; New = \"481\u00f72=240, 1\" },\n    @{ Old = \"698\u00f79=77, 5\"; New = \"689\u00f75=137, 4\" },\n    @{ Old = \"862\u00f78=107, 6\"; New = \"531\u00f73=177, 0\" },\n    @{ Old = \"104\u00f76=17, 2\"; New = \"803\u00f78=100, 3\" },\n    @{ Old = \"653\u00f73=217, 2\"; New = \"607\u00f73=202, 1\" },\n    @{ Old = \"335\u00f76=55, 5\"; New = \"635\u00f78=79, 3\" },\n    @{ Old = \"615\u00f77=87, 6\"; New = \"407\u00f73=135, 2\" },\n    @{ Old = \"618\u00f77=88, 2\"; New = \"305\u00f79=33, 8\" },\n    @{ Old = \"241\u00f74=60, 1\"; New = \"650\u00f77=92, 6\" },\n    @{ Old = \"502\u00f78=62, 6\"; New = \"370\u00f73=123, 1\" },\n    @{ Old = \"279\u00f77=39, 6\"; New = \"319\u00f78=39, 7\" },\n    @{ Old = \"101\u00f75=20, 1\"; New = \"311\u00f76=51, 5\" },\n    @{ Old = \"701\u00f78=87, 5\"; New = \"978\u00f78=122, 2\" },\n    @{ Old = \"703\u00f72=351, 1\"; New = \"286\u00f76=47, 4\" },\n    @{ Old = \"974\u00f77=139, 1\"; New = \"642\u00f77=91, 5\" },\n    @{ Old = \"956\u00f78=119, 4\"; New = \"333\u00f74=83, 1\" },\n    @{ Old = \"984\u00f75=196, 4\"; New = \"145\u00f78=18, 1\" },\n    @{ Old = \"206\u00f73=68, 2\"; New = \"896\u00f75=179, 1\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
